$wb = $excel.ActiveWorkbook
$wsMappings = $wb.Worksheets.Item(1)
$wsNomenclature = $wb.Worksheets.Item(2)

# Fill in the regex (column B) and case (column C) columns on the
# "nomenclature" sheet for the six channel rows.
$wsNomenclature.Range("B2").Value = "FSC\.H"
$wsNomenclature.Range("B3").Value = "SSC\.H"
$wsNomenclature.Range("B4").Value = "FL1\.H"
$wsNomenclature.Range("B5").Value = "FL2\.H"
$wsNomenclature.Range("B6").Value = "FL3\.H"
$wsNomenclature.Range("B7").Value = "FL4\.H"
$wsNomenclature.Range("C2:C7").Value = 0

# Update the saved selection / active sheet so the workbook re-opens with
# the "nomenclature" sheet active and its B7 cell selected, while the
# "mappings" sheet keeps A4 selected (but is no longer the active tab).
$null = $wsMappings.Range("A4").Select()
$null = $wsNomenclature.Activate()
$null = $wsNomenclature.Range("B7").Select()
